# Update countries & provincias Spain
#
# This script applies the data refresh represented in the diff:
#  - Updated case counts for several countries (Pakistan, Kazajistan,
#    Kirguistan, Angola).
#  - Liechtenstein / Butan swapped position in the underlying ranking,
#    so row 186 is now Butan and row 187 is now Liechtenstein (with
#    their respective updated numbers).
#  - Groenlandia / Islas Malvinas swapped position in the same way,
#    so row 209 is now Groenlandia and row 210 is now Islas Malvinas
#    (their figures were already identical, so only the names move).
#  - The "updated at" timestamp cell changed from 05:41 to 06:58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 06:58"

# --- Row 15: Pakistan ---
$ws.Range("B15").Value = 251625
$ws.Range("C15").Value = 2753
$ws.Range("D15").Value = 161917
$ws.Range("E15").Value = 84442
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 69
$ws.Range("H15").Value = 5266

# --- Row 34: Kazajistan ---
$ws.Range("B34").Value = 59899
$ws.Range("C34").Value = 1646
$ws.Range("D34").Value = 34190
$ws.Range("E34").Value = 25334
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 375

# --- Row 72: Kirguistan ---
$ws.Range("B72").Value = 11117
$ws.Range("C72").Value = 488
$ws.Range("D72").Value = 3460
$ws.Range("E72").Value = 7510
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 15
$ws.Range("H72").Value = 147

# --- Row 158: Angola ---
$ws.Range("B158").Value = 506
$ws.Range("C158").Value = 23
$ws.Range("D158").Value = 118
$ws.Range("E158").Value = 362
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 26

# --- Row 186: now Butan (was Liechtenstein) ---
$ws.Range("A186").Value = "Butan"
$ws.Range("B186").Value = 84
$ws.Range("C186").Value = 2
$ws.Range("D186").Value = 76
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# --- Row 187: now Liechtenstein (was Butan) ---
$ws.Range("A187").Value = "Liechtenstein"
$ws.Range("B187").Value = 84
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 81
$ws.Range("E187").Value = 2
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 1

# --- Row 209: now Groenlandia (was Islas Malvinas) ---
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# --- Row 210: now Islas Malvinas (was Groenlandia) ---
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
